$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextCell "D2" "67.667.81"
Set-TextCell "E2" "  +1.27%  "
Set-TextCell "D3" "2.621.05"
Set-TextCell "E3" "  +0.88%  "
Set-TextCell "E4" "  -0.13%  "
Set-TextCell "D5" "603.17"
Set-TextCell "E5" "  +1.72%  "
Set-TextCell "D6" "154.51"
Set-TextCell "E6" "  +0.60%  "
Set-TextCell "E7" "  +0.00%  "
Set-TextCell "E8" "  +1.59%  "
Set-TextCell "D9" "2.619.59"
Set-TextCell "E9" "  +0.92%  "
Set-TextCell "D10" "0.126"
Set-TextCell "E10" "  +9.97%  "
Set-TextCell "E11" "  +0.81%  "
Set-TextCell "E12" "  +1.05%  "
Set-TextCell "D13" "0.355"
Set-TextCell "E13" "  -0.48%  "
Set-TextCell "D14" "28.02"
Set-TextCell "E14" "  -0.23%  "
Set-TextCell "E15" "  +3.79%  "
Set-TextCell "D16" "3.095.43"
Set-TextCell "E16" "  +1.42%  "
Set-TextCell "D17" "67.546.01"
Set-TextCell "E17" "  +1.29%  "
Set-TextCell "D18" "2.621.00"
Set-TextCell "E18" "  +0.89%  "
Set-TextCell "D19" "11.30"
Set-TextCell "E19" "  +0.02%  "
Set-TextCell "D20" "363.98"
Set-TextCell "E20" "  +3.03%  "
Set-TextCell "E21" "  -2.78%  "
Set-TextCell "E22" "  -0.36%  "
Set-TextCell "D23" "2.13"
Set-TextCell "E23" "  +5.68%  "
Set-TextCell "D25" "70.14"
Set-TextCell "E25" "  +3.82%  "
Set-TextCell "D26" "10.12"
Set-TextCell "E26" "  -2.51%  "
Set-TextCell "E27" "  +3.25%  "
Set-TextCell "B28" "WrappedeETH"
Set-TextCell "C28" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextCell "D28" "2.745.13"
Set-TextCell "E28" "  +0.68%  "
Set-TextCell "B29" "Bittensor"
Set-TextCell "C29" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D29" "589.57"
Set-TextCell "E29" "  -0.12%  "
Set-TextCell "D30" "1.03"
Set-TextCell "E30" "  +3.12%  "
Set-TextCell "E31" "  -0.45%  "
Set-TextCell "E32" "  -0.51%  "
Set-TextCell "E33" "  +0.44%  "
Set-TextCell "E34" "  -2.45%  "
Set-TextCell "E35" "  +0.03%  "
Set-TextCell "E36" "  -0.81%  "
Set-TextCell "E37" "  -0.36%  "
Set-TextCell "E38" "  +1.01%  "
Set-TextCell "D39" "156.48"
Set-TextCell "E39" "  +2.38%  "
Set-TextCell "E40" "  +0.99%  "
Set-TextCell "D41" "5.44"
Set-TextCell "E41" "  +0.22%  "
Set-TextCell "E42" "  +3.03%  "
Set-TextCell "E43" "  +3.37%  "
Set-TextCell "D44" "41.15"
Set-TextCell "E44" "  -0.31%  "
Set-TextCell "B45" "WhiteBITCoin"
Set-TextCell "C45" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell "D45" "16.43"
Set-TextCell "E45" "  -0.16%  "
Set-TextCell "B46" "USDe"
Set-TextCell "C46" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D46" "0.999"
Set-TextCell "E46" "  -0.01%  "
Set-TextCell "D47" "157.13"
Set-TextCell "E47" "  +0.76%  "
Set-TextCell "D48" "0.0₆0288"
Set-TextCell "E48" "  -6.04%  "
Set-TextCell "E49" "  +0.39%  "
Set-TextCell "D50" "21.09"
Set-TextCell "E50" "  -0.26%  "
Set-TextCell "D51" "0.625"
Set-TextCell "E51" "  +1.08%  "
